# Weekly data refresh: a new sample ("Albahaca", Femacal de La Calera) is
# inserted as the new row 15, pushing the existing rows 15-101 down to
# 16-102 (dimension grows from A1:R101 to A1:R102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15, shifting rows 15..101 down to 16..102.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with this week's record.
$ws.Cells.Item(15, 1).Value  = 3
$ws.Cells.Item(15, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(15, 3).Value  = "Coquimbo"
$ws.Cells.Item(15, 4).Value  = 44550
$ws.Cells.Item(15, 5).Value  = 5
$ws.Cells.Item(15, 6).Value  = 100112052
$ws.Cells.Item(15, 7).Value  = "Albahaca"
$ws.Cells.Item(15, 8).Value  = "Sin especificar"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 140
$ws.Cells.Item(15, 11).Value = 4000
$ws.Cells.Item(15, 12).Value = 4500
$ws.Cells.Item(15, 13).Value = 4286
$ws.Cells.Item(15, 14).Value = "`$/docena de matas"
$ws.Cells.Item(15, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 16).Value = 714
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(15, 18).Value = "Hortaliza"
